$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.731.64'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '2.904.00'
$ws.Range("E3").Value = '  -2.53%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.32'
$ws.Range("E5").Value = '  -2.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.62'
$ws.Range("E6").Value = '  -5.37%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.546'
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("D9").Value = '2.911.36'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  -4.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.04'
$ws.Range("E11").Value = '  -1.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.358'
$ws.Range("E12").Value = '  -2.93%  '
$ws.Range("D13").Value = '3.408.53'
$ws.Range("E13").Value = '  -2.61%  '
$ws.Range("E14").Value = '  +3.28%  '
$ws.Range("D15").Value = '60.687.03'
$ws.Range("E15").Value = '  -1.73%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.55'
$ws.Range("E16").Value = '  -5.82%  '
$ws.Range("D17").Value = '2.908.05'
$ws.Range("E17").Value = '  -2.51%  '
$ws.Range("E18").Value = '  -4.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.95'
$ws.Range("E19").Value = '  -4.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  -3.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '351.57'
$ws.Range("E21").Value = '  -7.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.53'
$ws.Range("E22").Value = '  -2.76%  '
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.95'
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").Value = '  -4.25%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.177'
$ws.Range("E27").Value = '  -6.17%  '
$ws.Range("E28").Value = '  -0.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.84'
$ws.Range("E29").Value = '  -3.55%  '
$ws.Range("D30").Value = '0.0₃0857'
$ws.Range("E30").Value = '  -8.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.58'
$ws.Range("E33").Value = '  -4.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '152.38'
$ws.Range("E34").Value = '  -4.50%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.35'
$ws.Range("E35").Value = '  -4.71%  '
$ws.Range("E36").Value = '  -6.02%  '
$ws.Range("E37").Value = '  -7.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.19'
$ws.Range("E38").Value = '  -5.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.45'
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.46'
$ws.Range("E40").Value = '  -5.04%  '
$ws.Range("E41").Value = '  -4.34%  '
$ws.Range("D42").Value = '2.292.09'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.649'
$ws.Range("E43").Value = '  -3.26%  '
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.36'
$ws.Range("E45").Value = '  -7.56%  '
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.94'
$ws.Range("E47").Value = '  -3.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0237'
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.32'
$ws.Range("E49").Value = '  -0.94%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0916'
$ws.Range("E50").Value = '  -3.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.39'
